$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(6, 13).Value = 1.18
$ws.Cells.Item(6, 14).Value = 4.5
$ws.Cells.Item(6, 15).Value = 1.83
$ws.Cells.Item(6, 16).Value = 1.83
$ws.Cells.Item(6, 19).Value = 4
$ws.Cells.Item(6, 20).Value = 1.25
$ws.Cells.Item(6, 21).Value = 9
$ws.Cells.Item(6, 22).Value = 1.07
$ws.Cells.Item(6, 23).Value = 1.85
$ws.Cells.Item(6, 24).Value = 1.95
$ws.Cells.Item(6, 44).Value = 7
$ws.Cells.Item(6, 45).Value = 1.1

$ws.Cells.Item(7, 9).Value = 3.8
$ws.Cells.Item(7, 11).Value = 1.73
$ws.Cells.Item(7, 27).Value = 4.75
$ws.Cells.Item(7, 30).Value = 21
$ws.Cells.Item(7, 38).Value = 7
$ws.Cells.Item(7, 39).Value = 17
$ws.Cells.Item(7, 40).Value = 17
$ws.Cells.Item(7, 43).Value = 67

$ws.Cells.Item(8, 8).Value = 3.8
$ws.Cells.Item(8, 13).Value = 1.08
$ws.Cells.Item(8, 14).Value = 7.5
$ws.Cells.Item(8, 17).Value = 1.78
$ws.Cells.Item(8, 18).Value = 2.1
$ws.Cells.Item(8, 19).Value = 2.35
$ws.Cells.Item(8, 20).Value = 1.57
$ws.Cells.Item(8, 21).Value = 4.33
$ws.Cells.Item(8, 22).Value = 1.2
$ws.Cells.Item(8, 23).Value = 1.5
$ws.Cells.Item(8, 24).Value = 2.5
$ws.Cells.Item(8, 25).Value = 2.38
$ws.Cells.Item(8, 26).Value = 1.53
$ws.Cells.Item(8, 30).Value = 10
$ws.Cells.Item(8, 31).Value = 15
$ws.Cells.Item(8, 33).Value = 7.5
$ws.Cells.Item(8, 35).Value = 23
$ws.Cells.Item(8, 44).Value = 3.65
$ws.Cells.Item(8, 45).Value = 1.29

$ws.Cells.Item(11, 13).Value = 1.14
$ws.Cells.Item(11, 14).Value = 5.5
$ws.Cells.Item(11, 21).Value = 7
$ws.Cells.Item(11, 22).Value = 1.1

$ws.Cells.Item(15, 8).Value = 3.3
$ws.Cells.Item(15, 10).Value = 3
$ws.Cells.Item(15, 15).Value = 1.33
$ws.Cells.Item(15, 16).Value = 3.25
$ws.Cells.Item(15, 17).Value = 1.54
$ws.Cells.Item(15, 18).Value = 2.44
$ws.Cells.Item(15, 19).Value = 2.08
$ws.Cells.Item(15, 20).Value = 1.73
$ws.Cells.Item(15, 21).Value = 3.75
$ws.Cells.Item(15, 22).Value = 1.25
$ws.Cells.Item(15, 23).Value = 1.44
$ws.Cells.Item(15, 24).Value = 2.63
$ws.Cells.Item(15, 25).Value = 1.83
$ws.Cells.Item(15, 26).Value = 1.83
$ws.Cells.Item(15, 27).Value = 7.5
$ws.Cells.Item(15, 28).Value = 10
$ws.Cells.Item(15, 33).Value = 9
$ws.Cells.Item(15, 38).Value = 9.5
$ws.Cells.Item(15, 44).Value = 2.9
$ws.Cells.Item(15, 45).Value = 1.4

$ws.Cells.Item(16, 15).Value = 1.29
$ws.Cells.Item(16, 16).Value = 3.5
$ws.Cells.Item(16, 19).Value = 1.98
$ws.Cells.Item(16, 20).Value = 1.88
$ws.Cells.Item(16, 23).Value = 1.4

$ws.Cells.Item(17, 7).Value = 3.25
$ws.Cells.Item(17, 8).Value = 2.75
$ws.Cells.Item(17, 9).Value = 2.5
$ws.Cells.Item(17, 13).Value = 1.14
$ws.Cells.Item(17, 14).Value = 5.5
$ws.Cells.Item(17, 15).Value = 1.67
$ws.Cells.Item(17, 16).Value = 2.1
$ws.Cells.Item(17, 19).Value = 3.1
$ws.Cells.Item(17, 20).Value = 1.36
$ws.Cells.Item(17, 21).Value = 6.5
$ws.Cells.Item(17, 22).Value = 1.11
$ws.Cells.Item(17, 23).Value = 1.67
$ws.Cells.Item(17, 24).Value = 2.1
$ws.Cells.Item(17, 25).Value = 2.25
$ws.Cells.Item(17, 26).Value = 1.57
$ws.Cells.Item(17, 28).Value = 15
$ws.Cells.Item(17, 32).Value = 51
$ws.Cells.Item(17, 35).Value = 21
$ws.Cells.Item(17, 39).Value = 10
$ws.Cells.Item(17, 44).Value = 5
$ws.Cells.Item(17, 45).Value = 1.16

$ws.Cells.Item(18, 23).Value = 1.62
$ws.Cells.Item(18, 28).Value = 9

$ws.Cells.Item(20, 7).Value = 1.62
$ws.Cells.Item(20, 8).Value = 3.9
$ws.Cells.Item(20, 9).Value = 5
$ws.Cells.Item(20, 10).Value = 2.2
$ws.Cells.Item(20, 12).Value = 5.5
$ws.Cells.Item(20, 14).Value = 12
$ws.Cells.Item(20, 23).Value = 1.36
$ws.Cells.Item(20, 24).Value = 3
$ws.Cells.Item(20, 25).Value = 1.8
$ws.Cells.Item(20, 26).Value = 1.91
$ws.Cells.Item(20, 28).Value = 8
$ws.Cells.Item(20, 30).Value = 12
$ws.Cells.Item(20, 35).Value = 17
$ws.Cells.Item(20, 37).Value = 251
$ws.Cells.Item(20, 40).Value = 17

$ws.Cells.Item(34, 20).Value = 1.5

$ws.Cells.Item(46, 7).Value = 2.77
$ws.Cells.Item(46, 8).Value = 2.9
$ws.Cells.Item(46, 9).Value = 2.6
$ws.Cells.Item(46, 10).Value = 3.4
$ws.Cells.Item(46, 11).Value = 1.93
$ws.Cells.Item(46, 12).Value = 3.25
$ws.Cells.Item(46, 15).Value = 1.44
$ws.Cells.Item(46, 16).Value = 2.42
$ws.Cells.Item(46, 21).Value = 3.8
$ws.Cells.Item(46, 23).Value = 1.5
$ws.Cells.Item(46, 24).Value = 2.27
$ws.Cells.Item(46, 25).Value = 1.9
$ws.Cells.Item(46, 26).Value = 1.72
$ws.Cells.Item(46, 27).Value = 7.2
$ws.Cells.Item(46, 28).Value = 13
$ws.Cells.Item(46, 30).Value = 35
$ws.Cells.Item(46, 31).Value = 27
$ws.Cells.Item(46, 33).Value = 6.9
$ws.Cells.Item(46, 34).Value = 5.7
$ws.Cells.Item(46, 35).Value = 16
$ws.Cells.Item(46, 37).Value = 900
